$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Derrick White, PG,SG, Boston Celtics -> LaMelo Ball, PG,SG, Charlotte Hornets
$ws.Range("A4").Value = "LaMelo Ball"
$ws.Range("B4").Value = "PG,SG"
$ws.Range("C4").Value = "Charlotte Hornets"

# Row 5: LaMelo Ball, PG,SG, Charlotte Hornets -> Gradey Dick, SG,SF, Toronto Raptors
$ws.Range("A5").Value = "Gradey Dick"
$ws.Range("B5").Value = "SG,SF"
$ws.Range("C5").Value = "Toronto Raptors"

# Row 14: Grant Williams, PF,C, Charlotte Hornets -> Derrick White, PG,SG, Boston Celtics
$ws.Range("A14").Value = "Derrick White"
$ws.Range("B14").Value = "PG,SG"
$ws.Range("C14").Value = "Boston Celtics"

# Row 15: Brandon Miller, SG,SF, Charlotte Hornets -> Grant Williams, PF,C, Charlotte Hornets
$ws.Range("A15").Value = "Grant Williams"
$ws.Range("B15").Value = "PF,C"
$ws.Range("C15").Value = "Charlotte Hornets"

# Row 16: Gradey Dick, SG,SF, Toronto Raptors -> Brandon Miller, SG,SF, Charlotte Hornets
$ws.Range("A16").Value = "Brandon Miller"
$ws.Range("B16").Value = "SG,SF"
$ws.Range("C16").Value = "Charlotte Hornets"
